$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(2, 15.85831168482834, 12.63945213387992, 18.92867860913989),
    @(3, 17.69684784707468, 13.61691913763063, 21.47189480070147),
    @(4, 17.1670140878365, 14.64608114246996, 19.86450588533677),
    @(5, 21.40601637464777, 16.93890992690268, 25.6035948897917),
    @(6, 21.25182729152693, 18.51402150611681, 24.08685669652202),
    @(7, 14.24643923135528, 11.57817942179415, 17.06253846085109),
    @(8, 4.875452821812425, 3.207000433750151, 7.323750774061892),
    @(9, 19.59092723058472, 17.2761610890342, 22.48915157788146),
    @(10, 4.940534701879073, 3.188934070114838, 7.803381616884126),
    @(11, 12.88397627119304, 9.166611337785609, 17.32083232877174),
    @(12, 21.32986611255614, 17.57538731305797, 25.11510314202779),
    @(13, 7.03308902926521, 4.469250347012947, 11.06956196423835),
    @(14, 9.485787625537537, 6.298063163637429, 13.60023422940988),
    @(15, 25.48147809740647, 20.15155541886025, 31.43924131942448),
    @(16, 21.2368800065592, 18.45154832622705, 24.09408536614571),
    @(17, 21.32162833696087, 17.70554538489739, 24.94269668813383),
    @(18, 17.14531120486389, 14.64294758764848, 19.81835430858394),
    @(19, 13.46384327159959, 9.247965110643952, 18.55663906555142),
    @(20, 24.10587916593659, 20.8351078425097, 27.53498420171085),
    @(21, 24.21272169008476, 21.1163469436474, 28.02707874055412),
    @(22, 4.941334080942064, 3.278256580250292, 7.474817068607104),
    @(23, 26.95613202808381, 21.6745839321687, 32.26564750192062),
    @(24, 21.32323045398253, 18.5383844016844, 24.2234348390437),
    @(25, 18.07315072054334, 15.86753218119703, 20.51082379811587),
    @(26, 24.10097828733611, 21.08185181469875, 27.46322651452349),
    @(27, 7.086754965669518, 4.463349230929065, 10.88667272194598),
    @(28, 15.12546337189659, 12.59360931082749, 18.13349599513829),
    @(29, 28.78689838660346, 23.8202307687513, 34.75594095845128),
    @(30, 20.05876837892758, 17.11298132017111, 22.98934727440846),
    @(31, 9.063780866810335, 6.149079716916478, 12.31517642412827),
    @(32, 6.723486076341496, 4.539586460153668, 9.911779863023797),
    @(33, 9.488771294180921, 6.393383767618546, 13.60420496769811),
    @(34, 18.00185547221949, 13.26064402157137, 22.53464116178023),
    @(35, 22.61969692679389, 19.79891179836313, 25.6775063906045),
    @(36, 14.25672547837477, 11.57324553755865, 16.94387624989922),
    @(37, 26.1350998657001, 22.22205520575348, 30.06911346133122),
    @(38, 9.677531949858944, 6.172235935730229, 13.66440728744229),
    @(39, 21.93833621318485, 18.08336784859254, 26.22586288867161),
    @(40, 20.81363723742491, 17.73215804123688, 24.06318936933185),
    @(41, 13.32245204728365, 11.10787418942105, 15.46509347539225),
    @(42, 13.09597104300955, 9.304385506167764, 17.61078408310367),
    @(43, 20.0581177081499, 17.11072486066232, 22.98755217153743),
    @(44, 11.04869832632787, 8.411918766435258, 13.87533784038574),
    @(45, 13.32720116828266, 9.055937144890075, 17.5489841992224),
    @(46, 8.904165118766553, 5.970305076466435, 12.18406382753753),
    @(47, 5.018957560836838, 3.19089990607358, 7.745789201948318),
    @(48, 21.43234540378484, 17.57555277806394, 25.21822223929305),
    @(49, 18.95751915449611, 16.07753431737913, 21.82826545187725),
    @(50, 26.08419400448929, 22.2264855473252, 29.96942497224523),
    @(51, 17.68222642076326, 13.61358276844428, 21.44773332503449),
    @(52, 8.718864145421254, 6.220896679627077, 11.83549854912799),
    @(53, 21.97761300182407, 17.92826279956578, 26.0611810594617),
    @(54, 25.38235944533855, 20.02608353634214, 31.31003790780434),
    @(55, 24.12241941926847, 21.01877208660854, 27.80073367778109),
    @(56, 9.1752042687122, 6.482146440264168, 12.39424881163156),
    @(57, 6.71985316093577, 4.468085208783333, 9.801452457003192),
    @(58, 26.09866291007773, 22.22235707438722, 29.89696611861953),
    @(59, 12.05437658029511, 8.735008018149474, 15.52856705927967),
    @(60, 27.18943429493288, 21.92090683452057, 32.74171184343079),
    @(61, 9.655016244653604, 6.41656882845882, 14.15185844989305),
    @(62, 12.89872021790817, 9.187482350178581, 17.24474691744372),
    @(63, 9.736563933887577, 6.162881423430731, 14.55931158677468),
    @(64, 27.19421131838776, 21.94463219035836, 32.7359574131038),
    @(65, 6.579679562515837, 4.294559259120946, 9.421931006669848),
    @(66, 17.44695048844174, 13.4108228653142, 21.79868677428284),
    @(67, 16.75050609352498, 13.12568535069761, 20.30815283185558),
    @(68, 22.7028317956696, 19.4312358141118, 26.13950062840323),
    @(69, 4.981704021492909, 3.233522541080089, 7.680892900515082),
    @(70, 20.65026496877794, 17.85176502008645, 23.89752369307757),
    @(71, 15.08689766139918, 13.03989922131687, 17.15973979756303),
    @(72, 17.62019369232228, 15.06401229017496, 20.33597446836358),
)

foreach ($row in $values) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}
